# Scheduled market-data refresh: updates computed profit/price columns (H:N)
# for specific Leve rows across all eight crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2509.9
$ws.Range("I40").Value = 3080
$ws.Range("J40").Value = 1939.8
$ws.Range("K40").Value = 3080
$ws.Range("L40").Value = 1939.8
$ws.Range("M40").Value = -2905
$ws.Range("N40").Value = -2289.8
# Row 43
$ws.Range("H43").Value = 11159.429
$ws.Range("J43").Value = 3742.5557
$ws.Range("L43").Value = 3742.5557
$ws.Range("N43").Value = -3880.5557
# Row 69
$ws.Range("H69").Value = 8000
$ws.Range("J69").Value = 8000
$ws.Range("L69").Value = 24000
$ws.Range("N69").Value = -25748
# Row 70
$ws.Range("H70").Value = 252575
$ws.Range("J70").Value = 335633.34
$ws.Range("L70").Value = 1006900.02
$ws.Range("N70").Value = -1007440.02
# Row 72
$ws.Range("H72").Value = 8000
$ws.Range("J72").Value = 8000
$ws.Range("L72").Value = 72000
$ws.Range("N72").Value = -80736
# Row 73
$ws.Range("H73").Value = 252575
$ws.Range("J73").Value = 335633.34
$ws.Range("L73").Value = 1006900.02
$ws.Range("N73").Value = -1008772.02
# Row 94
$ws.Range("H94").Value = 660.3333
$ws.Range("I94").Value = 660.3333
$ws.Range("K94").Value = 660.3333
$ws.Range("M94").Value = -209.3333
# Row 96
$ws.Range("H96").Value = 1099.6
$ws.Range("I96").Value = 2499.5
$ws.Range("J96").Value = 166.33333
$ws.Range("K96").Value = 7498.5
$ws.Range("L96").Value = 498.99999
$ws.Range("M96").Value = -6125.5
$ws.Range("N96").Value = -3244.99999
# Row 98
$ws.Range("H98").Value = 3231.2222
$ws.Range("I98").Value = 3177.4666
$ws.Range("J98").Value = 3500
$ws.Range("K98").Value = 3177.4666
$ws.Range("L98").Value = 3500
$ws.Range("M98").Value = -1679.4666
$ws.Range("N98").Value = -6496
# Row 122
$ws.Range("H122").Value = 3231.2222
$ws.Range("I122").Value = 3177.4666
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 9532.399800000001
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -7082.399800000001
$ws.Range("N122").Value = -15400
$ws = $wb.Worksheets.Item("ARM")
# Row 104
$ws.Range("H104").Value = 60943.8
$ws.Range("J104").Value = 60943.8
$ws.Range("L104").Value = 60943.8
$ws.Range("N104").Value = -67931.8
# Row 122
$ws.Range("H122").Value = 1701.3158
$ws.Range("I122").Value = 1390.2222
$ws.Range("K122").Value = 4170.6666
$ws.Range("M122").Value = -1720.6666
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 748.0270400000001
$ws.Range("I107").Value = 656.375
$ws.Range("K107").Value = 656.375
$ws.Range("M107").Value = 1263.625
# Row 124
$ws.Range("H124").Value = 29999
$ws.Range("J124").Value = 29999
$ws.Range("L124").Value = 29999
$ws.Range("N124").Value = -39819
# Row 134
$ws.Range("H134").Value = 2213.8545
$ws.Range("I134").Value = 2212.4905
$ws.Range("J134").Value = 2250
$ws.Range("K134").Value = 6637.4715
$ws.Range("L134").Value = 6750
$ws.Range("M134").Value = -4102.4715
$ws.Range("N134").Value = -11820
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 227.5
$ws.Range("I22").Value = 227.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 227.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 122.5
$ws.Range("N22").ClearContents()
# Row 58
$ws.Range("H58").Value = 3350.8572
$ws.Range("I58").Value = 3110.8572
$ws.Range("J58").Value = 3830.8572
$ws.Range("K58").Value = 3110.8572
$ws.Range("L58").Value = 3830.8572
$ws.Range("M58").Value = -2907.8572
$ws.Range("N58").Value = -4236.8572
# Row 75
$ws.Range("H75").Value = 12500
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
# Row 78
$ws.Range("H78").Value = 12500
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
# Row 99
$ws.Range("H99").Value = 2600
$ws.Range("I99").Value = 2600
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2600
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1102
$ws.Range("N99").ClearContents()
# Row 107
$ws.Range("H107").Value = 828.3333
$ws.Range("I107").Value = 541.2222
$ws.Range("J107").Value = 1043.6666
$ws.Range("K107").Value = 541.2222
$ws.Range("L107").Value = 1043.6666
$ws.Range("M107").Value = 1378.7778
$ws.Range("N107").Value = -4883.6666
# Row 126
$ws.Range("H126").Value = 2600
$ws.Range("I126").Value = 2600
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7800
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5330
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 2841.5925
$ws.Range("I132").Value = 2756.348
$ws.Range("K132").Value = 8269.044
$ws.Range("M132").Value = -5739.044
# Row 136
$ws.Range("H136").Value = 3350.8572
$ws.Range("I136").Value = 3110.8572
$ws.Range("J136").Value = 3830.8572
$ws.Range("K136").Value = 9332.571599999999
$ws.Range("L136").Value = 11492.5716
$ws.Range("M136").Value = -6782.571599999999
$ws.Range("N136").Value = -16592.5716
$ws = $wb.Worksheets.Item("CUL")
# Row 21
$ws.Range("H21").Value = 290.6
$ws.Range("I21").Value = 218.16667
$ws.Range("J21").Value = 399.25
$ws.Range("K21").Value = 654.50001
$ws.Range("L21").Value = 1197.75
$ws.Range("M21").Value = -481.50001
$ws.Range("N21").Value = -1543.75
# Row 107
$ws.Range("H107").Value = 682
$ws.Range("I107").Value = 802.8570999999999
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 2408.5713
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = -488.5712999999996
$ws.Range("N107").Value = -5040
$ws = $wb.Worksheets.Item("GSM")
# Row 18
$ws.Range("H18").Value = 27798022
$ws.Range("J18").Value = 39989.5
$ws.Range("L18").Value = 39989.5
$ws.Range("N18").Value = -40575.5
# Row 20
$ws.Range("H20").Value = 16000
$ws.Range("J20").Value = 16000
$ws.Range("L20").Value = 16000
$ws.Range("N20").Value = -16490
# Row 43
$ws.Range("H43").Value = 4904.25
$ws.Range("I43").Value = 4617
$ws.Range("J43").Value = 5000
$ws.Range("K43").Value = 4617
$ws.Range("L43").Value = 5000
$ws.Range("M43").Value = -4466
$ws.Range("N43").Value = -5302
# Row 46
$ws.Range("H46").Value = 22666.334
$ws.Range("I46").Value = 8999.5
$ws.Range("J46").Value = 50000
$ws.Range("K46").Value = 8999.5
$ws.Range("L46").Value = 50000
$ws.Range("M46").Value = -8843.5
$ws.Range("N46").Value = -50312
# Row 80
$ws.Range("H80").Value = 2999.5
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 3499
$ws.Range("K80").Value = 2500
$ws.Range("L80").Value = 3499
$ws.Range("M80").Value = -1502
$ws.Range("N80").Value = -5495
# Row 83
$ws.Range("H83").Value = 2999.5
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 3499
$ws.Range("K83").Value = 12500
$ws.Range("L83").Value = 17495
$ws.Range("M83").Value = -7508
$ws.Range("N83").Value = -27479
# Row 126
$ws.Range("H126").Value = 16416.75
$ws.Range("I126").Value = 19304.23
$ws.Range("J126").Value = 3904.3333
$ws.Range("K126").Value = 57912.69
$ws.Range("L126").Value = 11712.9999
$ws.Range("M126").Value = -55442.69
$ws.Range("N126").Value = -16652.9999
# Row 132
$ws.Range("H132").Value = 2983.1428
$ws.Range("I132").Value = 2782.3
$ws.Range("K132").Value = 8346.900000000001
$ws.Range("M132").Value = -5816.900000000001
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2319.125
$ws.Range("I16").Value = 2455.1904
$ws.Range("J16").Value = 1366.6666
$ws.Range("K16").Value = 2455.1904
$ws.Range("L16").Value = 1366.6666
$ws.Range("M16").Value = -2285.1904
$ws.Range("N16").Value = -1706.6666
# Row 46
$ws.Range("H46").Value = 1000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -1376
# Row 82
$ws.Range("H82").Value = 3354.25
$ws.Range("I82").Value = 3838.5
$ws.Range("J82").Value = 1901.5
$ws.Range("K82").Value = 3838.5
$ws.Range("L82").Value = 1901.5
$ws.Range("M82").Value = -3477.5
$ws.Range("N82").Value = -2623.5
# Row 85
$ws.Range("H85").Value = 3354.25
$ws.Range("I85").Value = 3838.5
$ws.Range("J85").Value = 1901.5
$ws.Range("K85").Value = 3838.5
$ws.Range("L85").Value = 1901.5
$ws.Range("M85").Value = -2590.5
$ws.Range("N85").Value = -4397.5
# Row 100
$ws.Range("H100").Value = 7999.5557
$ws.Range("I100").Value = 2699.4
$ws.Range("K100").Value = 2699.4
$ws.Range("M100").Value = -2158.4
# Row 132
$ws.Range("H132").Value = 4634.1113
$ws.Range("I132").Value = 4377.294
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 13131.882
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -10601.882
$ws.Range("N132").Value = -32060
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2635.8333
$ws.Range("I122").Value = 2326
$ws.Range("K122").Value = 6978
$ws.Range("M122").Value = -4528
# Row 132
$ws.Range("H132").Value = 2194.08
$ws.Range("I132").Value = 2242.35
$ws.Range("K132").Value = 6727.049999999999
$ws.Range("M132").Value = -4197.049999999999
# Row 136
$ws.Range("H136").Value = 2519.5833
$ws.Range("I136").Value = 2481.0715
$ws.Range("J136").Value = 2654.375
$ws.Range("K136").Value = 7443.2145
$ws.Range("L136").Value = 7963.125
$ws.Range("M136").Value = -4893.2145
$ws.Range("N136").Value = -13063.125
